# The deck's theme color scheme (applied via the slide master, theme2.xml)
# is switched from the "Integral" / "Red Violet" palette back to the
# stock Office theme's "Office" colour scheme.
#
# PowerPoint's ColorScheme.Item(n).RGB uses the usual COM BGR-packed
# integer (&H00BBGGRR), so build it from R/G/B bytes explicitly.
function ToRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

# Index -> theme colour slot (document order inside <a:clrScheme>):
#  1 dk1       2 lt1       3 dk2        4 lt2
#  5 accent1   6 accent2   7 accent3    8 accent4
#  9 accent5  10 accent6  11 hlink     12 folHlink
$cs.Item(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1      000000
$cs.Item(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$cs.Item(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2      44546A
$cs.Item(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$cs.Item(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$cs.Item(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2  ED7D31
$cs.Item(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$cs.Item(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4  FFC000
$cs.Item(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5  4472C4
$cs.Item(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6  70AD47
$cs.Item(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink    0563C1
$cs.Item(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink 954F72
